# Refactored projection output: recomputed Monte-Carlo style results and
# applied a Currency number format to the data table, plus moved the
# selected cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Summary")

# --- Updated projection results (rows 2-9, columns A:E) -------------------
$data = @(
    @(1000, 413274.62897982268, 30553.772961806899, 339029.13926513313, 25064.73570398978),
    @(1500, 523563.0729180249,  37774.158789359877, 429504.0767457136,  30987.96694206553),
    @(2000, 622061.2902191258,  41167.014828503197, 510306.92184940551, 33771.29062549791),
    @(2500, 737269.80480622686, 49659.178289565563, 604818.03092207585, 40737.822483042313),
    @(3000, 845897.35236389283, 56314.778021926068, 693930.45488061442, 46197.732408133808),
    @(3500, 951852.76737273403, 61534.393198757527, 780850.79944568395, 50479.634844450047),
    @(4000, 1060396.689416847,  68923.707077803323, 869894.62135635363, 56541.445922369159),
    @(5000, 1279712.4929762781, 81115.775720868696, 1049809.967942086,  66543.188705669018)
)

$row = 2
foreach ($r in $data) {
    $ws.Cells.Item($row, 1).Value = $r[0]
    $ws.Cells.Item($row, 2).Value = $r[1]
    $ws.Cells.Item($row, 3).Value = $r[2]
    $ws.Cells.Item($row, 4).Value = $r[3]
    $ws.Cells.Item($row, 5).Value = $r[4]
    $row = $row + 1
}

# --- Apply the built-in "Currency" cell style to the data table -----------
$ws.Range("A2:E9").Style = "Currency"

# --- Update the active selection shown when the sheet is reopened ---------
$ws.Range("D7").Select()

Write-Output "applied projection refresh + currency formatting"
